$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Style = 'Normal'
$ws.Range("D2").Value = "'59.449.75"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -6.80%  '
$ws.Range("D3").Style = 'Normal'
$ws.Range("D3").Value = "'3.306.55"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -3.81%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Style = 'Normal'
$ws.Range("D5").Value = "'557.73"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -4.19%  '
$ws.Range("D6").Style = 'Normal'
$ws.Range("D6").Value = "'127.64"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -2.25%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Style = 'Normal'
$ws.Range("D8").Value = "'3.306.52"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -3.79%  '
$ws.Range("D9").Style = 'Normal'
$ws.Range("D9").Value = "'0.468"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  -2.43%  '
$ws.Range("D10").Style = 'Normal'
$ws.Range("D10").Value = "'7.33"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -4.23%  '
$ws.Range("E11").Value = '  -6.26%  '
$ws.Range("D12").Style = 'Normal'
$ws.Range("D12").Value = "'0.370"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -3.49%  '
$ws.Range("D13").Style = 'Normal'
$ws.Range("D13").Value = "'3.863.66"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -4.05%  '
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").Style = 'Normal'
$ws.Range("D15").Value = "'3.288.34"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -4.40%  '
$ws.Range("D16").Style = 'Normal'
$ws.Range("D16").Value = "'0.0000166"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -6.22%  '
$ws.Range("D17").Style = 'Normal'
$ws.Range("D17").Value = "'24.15"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -3.96%  '
$ws.Range("D18").Style = 'Normal'
$ws.Range("D18").Value = "'59.621.47"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -6.45%  '
$ws.Range("D19").Style = 'Normal'
$ws.Range("D19").Value = "'5.61"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("D20").Style = 'Normal'
$ws.Range("D20").Value = "'13.23"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").Style = 'Normal'
$ws.Range("D21").Value = "'8.90"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -9.95%  '
$ws.Range("D22").Style = 'Normal'
$ws.Range("D22").Value = "'349.95"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -8.77%  '
$ws.Range("D23").Style = 'Normal'
$ws.Range("D23").Value = "'0.552"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -1.92%  '
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").Style = 'Normal'
$ws.Range("D25").Value = "'3.433.90"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -3.99%  '
$ws.Range("D26").Style = 'Normal'
$ws.Range("D26").Value = "'68.48"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -7.31%  '
$ws.Range("D27").Style = 'Normal'
$ws.Range("D27").Value = "'0.0000110"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").Style = 'Normal'
$ws.Range("D28").Value = "'0.995"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").Style = 'Normal'
$ws.Range("D29").Value = "'7.31"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +3.56%  '
$ws.Range("D30").Style = 'Normal'
$ws.Range("D30").Value = "'1.47"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +3.77%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Style = 'Normal'
$ws.Range("D31").Value = "'7.78"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").Style = 'Normal'
$ws.Range("D32").Value = "'0.151"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -2.25%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Style = 'Normal'
$ws.Range("D33").Value = "'2.08"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -5.86%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Style = 'Normal'
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").Style = 'Normal'
$ws.Range("D35").Value = "'3.328.15"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -3.99%  '
$ws.Range("D36").Style = 'Normal'
$ws.Range("D36").Value = "'22.59"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("D37").Style = 'Normal'
$ws.Range("D37").Value = "'5.30"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +2.07%  '
$ws.Range("D38").Style = 'Normal'
$ws.Range("D38").Value = "'6.74"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("D40").Style = 'Normal'
$ws.Range("D40").Value = "'157.06"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -3.82%  '
$ws.Range("D41").Style = 'Normal'
$ws.Range("D41").Value = "'0.0746"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -3.52%  '
$ws.Range("D42").Style = 'Normal'
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").Style = 'Normal'
$ws.Range("D43").Value = "'40.57"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").Style = 'Normal'
$ws.Range("D44").Value = "'1.18"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +5.76%  '
$ws.Range("D45").Style = 'Normal'
$ws.Range("D45").Value = "'4.28"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Style = 'Normal'
$ws.Range("D46").Value = "'0.739"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -6.55%  '
$ws.Range("D47").Style = 'Normal'
$ws.Range("D47").Value = "'22.92"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("D48").Style = 'Normal'
$ws.Range("D48").Value = "'1.53"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -5.01%  '
$ws.Range("D49").Style = 'Normal'
$ws.Range("D49").Value = "'6.70"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").Style = 'Normal'
$ws.Range("D50").Value = "'2.41"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +16.24%  '
$ws.Range("D51").Style = 'Normal'
$ws.Range("D51").Value = "'21.74"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +6.90%  '
